function Set-CellText($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "67.709.45"
Set-CellText $ws "E2" "  +1.37%  "
Set-CellText $ws "D3" "2.492.62"
Set-CellText $ws "E3" "  +1.60%  "
Set-CellText $ws "E4" "  -0.03%  "
Set-CellText $ws "D5" "586.63"
Set-CellText $ws "E5" "  +1.06%  "
Set-CellText $ws "D6" "176.81"
Set-CellText $ws "E6" "  +5.41%  "
Set-CellText $ws "D8" "0.516"
Set-CellText $ws "E8" "  +1.32%  "
Set-CellText $ws "E9" "  +4.63%  "
Set-CellText $ws "E10" "  +0.70%  "
Set-CellText $ws "E11" "  +3.99%  "
Set-CellText $ws "E12" "  +1.60%  "
Set-CellText $ws "D13" "2.931.31"
Set-CellText $ws "E13" "  +1.04%  "
Set-CellText $ws "D14" "25.74"
Set-CellText $ws "E14" "  +2.41%  "
Set-CellText $ws "D15" "67.629.58"
Set-CellText $ws "E15" "  +1.39%  "
Set-CellText $ws "E16" "  +2.26%  "
Set-CellText $ws "D17" "2.484.37"
Set-CellText $ws "E17" "  +2.85%  "
Set-CellText $ws "D18" "11.06"
Set-CellText $ws "E18" "  +1.38%  "
Set-CellText $ws "E19" "  +0.64%  "
Set-CellText $ws "D20" "351.74"
Set-CellText $ws "E20" "  +0.86%  "
Set-CellText $ws "E21" "  +1.73%  "
Set-CellText $ws "E22" "  +0.22%  "
Set-CellText $ws "D23" "70.71"
Set-CellText $ws "E23" "  +3.20%  "
Set-CellText $ws "D24" "4.24"
Set-CellText $ws "E24" "  +1.43%  "
Set-CellText $ws "E25" "  -0.19%  "
Set-CellText $ws "E26" "  +2.27%  "
Set-CellText $ws "D27" "2.620.22"
Set-CellText $ws "E27" "  +1.57%  "
Set-CellText $ws "D28" "0.998"
Set-CellText $ws "E28" "  -0.59%  "
Set-CellText $ws "D29" "0.0₃0910"
Set-CellText $ws "E29" "  +2.52%  "
Set-CellText $ws "D30" "515.41"
Set-CellText $ws "E30" "  +2.23%  "
Set-CellText $ws "D31" "7.87"
Set-CellText $ws "E31" "  +4.02%  "
Set-CellText $ws "E32" "  +3.60%  "
Set-CellText $ws "E33" "  +1.83%  "
Set-CellText $ws "E34" "  +0.01%  "
Set-CellText $ws "E35" "  +7.67%  "
Set-CellText $ws "D36" "160.89"
Set-CellText $ws "E36" "  +1.74%  "
Set-CellText $ws "E37" "  +0.37%  "
Set-CellText $ws "D38" "18.38"
Set-CellText $ws "E38" "  +1.35%  "
Set-CellText $ws "E39" "  +2.37%  "
Set-CellText $ws "E40" "  +0.10%  "
Set-CellText $ws "D41" "1.73"
Set-CellText $ws "E41" "  +4.22%  "
Set-CellText $ws "B42" "RenderToken"
Set-CellText $ws "C42" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-CellText $ws "D42" "4.89"
Set-CellText $ws "E42" "  +3.37%  "
Set-CellText $ws "B43" "PolygonEcosystemToken"
Set-CellText $ws "C43" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-CellText $ws "D43" "0.330"
Set-CellText $ws "E43" "  +2.47%  "
Set-CellText $ws "E44" "  +4.26%  "
Set-CellText $ws "D45" "144.07"
Set-CellText $ws "E45" "  +3.02%  "
Set-CellText $ws "E46" "  +3.22%  "
Set-CellText $ws "E47" "  +2.23%  "
Set-CellText $ws "D48" "0.0747"
Set-CellText $ws "E48" "  +2.92%  "
Set-CellText $ws "E49" "  +1.79%  "
Set-CellText $ws "E50" "  +1.71%  "
Set-CellText $ws "E51" "  +1.65%  "
